$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.096.73"
$ws.Cells.Item(2, 5).Value = "  +3.10%  "
$ws.Cells.Item(3, 4).Value = "2.311.60"
$ws.Cells.Item(3, 5).Value = "  +2.88%  "
$ws.Cells.Item(4, 5).Value = "  -0.09%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "310.08"
$ws.Cells.Item(5, 5).Value = "  +2.12%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "100.85"
$ws.Cells.Item(6, 5).Value = "  +6.49%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.536"
$ws.Cells.Item(7, 5).Value = "  +2.62%  "
$ws.Cells.Item(8, 5).Value = "  -0.05%  "
$ws.Cells.Item(9, 5).Value = "  +8.10%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "36.06"
$ws.Cells.Item(10, 5).Value = "  +4.60%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0815"
$ws.Cells.Item(11, 5).Value = "  +3.92%  "
$ws.Cells.Item(12, 5).Value = "  +0.99%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "7.20"
$ws.Cells.Item(13, 5).Value = "  +6.54%  "
$ws.Cells.Item(14, 4).Value = "2.667.02"
$ws.Cells.Item(14, 5).Value = "  +2.75%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "14.97"
$ws.Cells.Item(15, 5).Value = "  +4.48%  "
$ws.Cells.Item(16, 4).Value = "2.323.14"
$ws.Cells.Item(16, 5).Value = "  +3.02%  "
$ws.Cells.Item(17, 5).Value = "  +3.73%  "
$ws.Cells.Item(18, 4).Value = "43.023.31"
$ws.Cells.Item(18, 5).Value = "  +3.20%  "
$ws.Cells.Item(19, 5).Value = "  +3.08%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0919"
$ws.Cells.Item(20, 5).Value = "  +2.50%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.12"
$ws.Cells.Item(21, 5).Value = "  +3.53%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "68.44"
$ws.Cells.Item(22, 5).Value = "  +0.46%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "240.84"
$ws.Cells.Item(23, 5).Value = "  +2.21%  "
$ws.Cells.Item(24, 5).Value = "  +6.56%  "
$ws.Cells.Item(25, 5).Value = "  +3.50%  "
$ws.Cells.Item(26, 5).Value = "  +0.08%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "24.71"
$ws.Cells.Item(27, 5).Value = "  +5.64%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "37.33"
$ws.Cells.Item(28, 5).Value = "  +3.62%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.65"
$ws.Cells.Item(29, 5).Value = "  +3.15%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.11"
$ws.Cells.Item(30, 5).Value = "  -0.22%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "167.22"
$ws.Cells.Item(31, 5).Value = "  +4.70%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "5.33"
$ws.Cells.Item(32, 5).Value = "  +3.61%  "
$ws.Cells.Item(33, 5).Value = "  -0.04%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.17"
$ws.Cells.Item(34, 5).Value = "  +0.59%  "
$ws.Cells.Item(35, 5).Value = "  +7.49%  "
$ws.Cells.Item(36, 5).Value = "  +2.04%  "
$ws.Cells.Item(37, 5).Value = "  +3.45%  "
$ws.Cells.Item(38, 5).Value = "  +0.83%  "
$ws.Cells.Item(39, 5).Value = "  +3.35%  "
$ws.Cells.Item(40, 5).Value = "  +2.38%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "4.29"
$ws.Cells.Item(41, 5).Value = "  +8.81%  "
$ws.Cells.Item(42, 5).Value = "  +1.98%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "19.61"
$ws.Cells.Item(43, 5).Value = "  +6.20%  "
$ws.Cells.Item(44, 5).Value = "  +3.95%  "
$ws.Cells.Item(45, 4).Value = "1.976.50"
$ws.Cells.Item(45, 5).Value = "  +1.15%  "
$ws.Cells.Item(46, 5).Value = "  +4.94%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.82"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.97"
$ws.Cells.Item(48, 5).Value = "  +18.92%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "55.59"
$ws.Cells.Item(49, 5).Value = "  +5.68%  "
$ws.Cells.Item(50, 4).Value = "2.536.79"
$ws.Cells.Item(50, 5).Value = "  +2.65%  "
$ws.Cells.Item(51, 5).Value = "  +4.89%  "
